$d = $word.ActiveDocument

# The first paragraph currently reads "Version 1." as:
#   [spellStart]<r>Version</r>[spellEnd] <r> 1.</r> [bookmarkStart/_GoBack][bookmarkEnd]
# It needs to become "Version 2." but re-split across runs as:
#   [spellStart]<r>Versi</r><r>on</r>[spellEnd] <r> 2</r> [bookmarkStart/_GoBack][bookmarkEnd] <r>.</r>
#
# We rebuild the paragraph's content (as raw WordprocessingML) via InsertXML
# over the full paragraph range (including its end-of-paragraph mark), which
# lets us control the run boundaries precisely. InsertXML on a range that
# spans a paragraph mark inserts the new content *and* leaves the original
# (now empty) paragraph mark behind as a trailing empty paragraph, so we
# merge that leftover paragraph back afterwards.

$p = $d.Paragraphs.Item(1)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="079893AC" w14:textId="77777777" w:rsidR="00F87116" w:rsidRDefault="001C5B2D"><w:proofErr w:type="spellStart"/><w:r><w:t>Versi</w:t></w:r><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$p.Range.InsertXML($xml)

# InsertXML left the paragraph mark from the old paragraph as a new, empty
# trailing paragraph. Delete the boundary between the freshly inserted
# paragraph and that leftover empty one so the document is back to a single
# paragraph again.
$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
[void]$d.Range($p1.Range.End - 1, $p2.Range.End).Delete()
